# v0.7.4 testing now includes BC trials as well as CB trials.
# Append the four "BC" trial rows (6-9) under the existing "CB" trial rows
# (2-5) on the testing-block worksheet, mirroring the existing layout:
#   A: sample_category   B: required_response
#   C: target_category_left   D: target_category_right
#   E: trial_description

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 6; A = "b1"; B = "'left";  C = "c1"; D = "c2"; E = "bc" },
    @{ Row = 7; A = "b1"; B = "right";  C = "c2"; D = "c1"; E = "bc" },
    @{ Row = 8; A = "b2"; B = "'left";  C = "c2"; D = "c1"; E = "bc" },
    @{ Row = 9; A = "b2"; B = "right";  C = "c1"; D = "c2"; E = "bc" }
)

foreach ($r in $newRows) {
    # Match the vertical-center alignment style used by the existing
    # "CB" trial rows (2-5) so the new rows share the same cell style.
    # Applying this before setting values avoids minting an extra,
    # unused intermediate style entry.
    $ws.Range("A" + $r.Row + ":E" + $r.Row).VerticalAlignment = -4108

    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
}

$ws.Range("E9").Select() | Out-Null
